$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.787.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "'2.282.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'251.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "

$ws.Range("D6").Value = "'0.645"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.51%  "

$ws.Range("D7").Value = "'74.74"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.44%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.645"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.94%  "

$ws.Range("D10").Value = "'39.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.58%  "

$ws.Range("D11").Value = "'0.0981"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.86%  "

$ws.Range("D12").Value = "'7.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("D13").Value = "'0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.50%  "

$ws.Range("D14").Value = "'2.628.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.84%  "

$ws.Range("D15").Value = "'15.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.54%  "

$ws.Range("D16").Value = "'0.872"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.64%  "

$ws.Range("D17").Value = "'2.269.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("D18").Value = "'42.695.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("E19").Value = "  +1.88%  "

$ws.Range("D20").Value = "'6.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.62%  "

$ws.Range("D21").Value = "'72.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.50%  "

$ws.Range("D22").Value = "'237.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.80%  "

$ws.Range("E23").Value = "  +5.59%  "

$ws.Range("D24").Value = "'3.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.33%  "

$ws.Range("D26").Value = "'11.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.69%  "

$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("E28").Value = "  +2.73%  "

$ws.Range("D29").Value = "'167.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("D30").Value = "'21.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("D31").Value = "'0.0890"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.74%  "

$ws.Range("D32").Value = "'6.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "

$ws.Range("D33").Value = "'0.127"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("D34").Value = "'31.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("D35").Value = "'0.128"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.18%  "

$ws.Range("D36").Value = "'4.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.11%  "

$ws.Range("D37").Value = "'4.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.93%  "

$ws.Range("E38").Value = "  -4.01%  "

$ws.Range("D39").Value = "'13.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.25%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").Value = "'5.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.29%  "

$ws.Range("D42").Value = "'0.210"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.79%  "

$ws.Range("D43").Value = "'9.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "'61.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.72%  "

$ws.Range("D45").Value = "'4.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "

$ws.Range("E47").Value = "  -1.11%  "

$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("E50").Value = "  -1.34%  "

$ws.Range("E51").Value = "  -0.72%  "
